# Add a new "classes" worksheet (after "components") that documents which
# fields/slots are used by the new workflowStepSet base class and its
# subclasses (featuresSet, featuresGroupsSet, MSPeakListsSet, formulasSet,
# compoundsSet, componentsSet).
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "classes"

# Fill data in the same order the original author typed it, so the
# shared-strings table indices line up with the target workbook.
$ws.Range("A2").Value = "featuresSet"
$ws.Range("B1").Value = "adducts"
$ws.Range("C1").Value = "setObjects"
$ws.Range("D1").Value = "ionizedXXX"
$ws.Range("A3").Value = "featuresGroupsSet"
$ws.Range("A4").Value = "MSPeakListsSet"
$ws.Range("A5").Value = "formulasSet"
$ws.Range("A6").Value = "compoundsSet"
$ws.Range("A7").Value = "componentsSet"
$ws.Range("E1").Value = "setThreshold"
$ws.Range("F1").Value = "origFGNames"
$ws.Range("G1").Value = "groupAlgorithm"
$ws.Range("H1").Value = "analysisInfo"

# X marks
$ws.Range("B2").Value = "X"
$ws.Range("C2").Value = "X"
$ws.Range("D2").Value = "X"

$ws.Range("G3").Value = "X"

$ws.Range("B4").Value = "X"
$ws.Range("C4").Value = "X"
$ws.Range("H4").Value = "X"

$ws.Range("B5").Value = "X"
$ws.Range("C5").Value = "X"
$ws.Range("E5").Value = "X"
$ws.Range("F5").Value = "X"

$ws.Range("B6").Value = "X"
$ws.Range("C6").Value = "X"
$ws.Range("E6").Value = "X"
$ws.Range("F6").Value = "X"

$ws.Range("C7").Value = "X"

# Center-align the X cells like the style used
$ws.Range("B2:I7").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B8:D10").HorizontalAlignment = -4108  # xlCenter

# Match the (auto-fit) column widths of the original workbook as closely
# as this environment's column-width model allows.
$ws.Columns.Item(1).ColumnWidth = 17.022135416666668
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666
$ws.Columns.Item(5).ColumnWidth = 11.736979166666666
$ws.Columns.Item(6).ColumnWidth = 12.022135416666666
$ws.Columns.Item(7).ColumnWidth = 14.307291666666666
$ws.Columns.Item(8).ColumnWidth = 10.736979166666666

$ws.Range("C3").Select() | Out-Null
